{"js": "// Load all paragraphs in the body so we can locate the anchor points by text.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1) The \"_GoBack\" bookmark currently sits at the end of the \"Test XGboost\"\n//    paragraph. It is being moved, so remove it from its old location first.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Find the heading paragraph \"Classification problem\"; the two new\n//    paragraphs are inserted right after it (i.e. right before the blank\n//    paragraph that currently follows the heading).\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Classification problem\") {\n    anchor = paragraphs.items[i + 1];\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error('Could not find the \"Classification problem\" heading.');\n}\n\n// 3) Insert the new \"Standardize/Normalize data\" paragraph, then the new\n//    blank paragraph above it, both placed before the pre-existing blank\n//    paragraph (`anchor`) so the final order is:\n//      Classification problem\n//      <blank>\n//      Standardize/Normalize data\n//      <pre-existing blank paragraph>\n//      Models : ...\nconst textParagraph = anchor.insertParagraph(\"Standardize/Normalize data\", \"Before\");\nawait context.sync();\n\nconst blankParagraph = textParagraph.insertParagraph(\"\", \"Before\");\nawait context.sync();\n\n// 4) Re-create the \"_GoBack\" bookmark around the newly inserted text.\nconst newTextRange = textParagraph.getRange(\"Content\");\nnewTextRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Helper: find the 1-based index of the paragraph that contains document\n# position $pos. We avoid relying on Paragraph.Next()/.Previous() chains\n# across mutations since those can get confused once the document has been\n# edited; plain integer indices into $d.Paragraphs stay reliable instead.\nfunction Get-ParagraphIndexAtPosition($doc, $pos) {\n    $paras = $doc.Paragraphs\n    for ($i = 1; $i -le $paras.Count; $i++) {\n        $p = $paras.Item($i)\n        if ($p.Range.Start -le $pos -and $pos -le $p.Range.End) {\n            return $i\n        }\n    }\n    return -1\n}\n\n$d = $word.ActiveDocument\n\n# 1) The \"_GoBack\" bookmark currently wraps (the end of) the \"Test XGboost\"\n#    paragraph. It is being relocated, so delete it from there first.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Locate the \"Classification problem\" heading paragraph; the two new\n#    paragraphs go right after it (i.e. right before the blank paragraph\n#    that currently follows the heading).\n$findRange = $d.Content\n$null = $findRange.Find.Execute(\"Classification problem\")\n$headingIdx = Get-ParagraphIndexAtPosition $d $findRange.Start\n$anchorIdx = $headingIdx + 1\n\n# 3) Insert a new blank paragraph right before the anchor paragraph.\n$d.Paragraphs.Item($anchorIdx).Range.InsertParagraphBefore()\n$anchorIdx = $anchorIdx + 1\n\n# 4) Insert a new paragraph (for the \"Standardize/Normalize data\" text)\n#    right before the anchor paragraph, i.e. right after the blank one.\n$d.Paragraphs.Item($anchorIdx).Range.InsertParagraphBefore()\n$textIdx = $anchorIdx\n$anchorIdx = $anchorIdx + 1\n\n$d.Paragraphs.Item($textIdx).Range.Text = \"Standardize/Normalize data\"\n\n# 5) Re-create the \"_GoBack\" bookmark around the new text (excluding the\n#    trailing paragraph mark, so the bookmark stays inside this paragraph).\n$bmRange = $d.Paragraphs.Item($textIdx).Range.Duplicate\n$bmRange.MoveEnd(1, -1)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
